$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column sometimes holds plain-looking numeric strings (e.g. "1.005",
# "16.50"). A leading apostrophe forces Excel to store the literal text (matching
# quotePrefix semantics) instead of silently coercing it to a Number, which would
# also drop significant trailing zeros (16.50 -> 16.5).

$ws.Range("D2").Value = "27.493.49"
$ws.Range("E2").Value = "  +5.48%  "

$ws.Range("D3").Value = "1.725.36"
$ws.Range("E3").Value = "  +4.78%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Value = "'225.56"
$ws.Range("E5").Value = "  +3.28%  "

$ws.Range("D6").Value = "'0.5348"
$ws.Range("E6").Value = "  +3.04%  "

$ws.Range("D7").Value = "'1.005"
$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").Value = "'0.2664"
$ws.Range("E8").Value = "  +1.46%  "

$ws.Range("E9").Value = "  +4.58%  "

$ws.Range("D10").Value = "'21.58"
$ws.Range("E10").Value = "  +6.61%  "

$ws.Range("D11").Value = "'0.07687"
$ws.Range("E11").Value = "  +0.01%  "

$ws.Range("D12").Value = "'4.602"
$ws.Range("E12").Value = "  +0.28%  "

$ws.Range("D13").Value = "1.725.57"
$ws.Range("E13").Value = "  +4.46%  "

$ws.Range("D14").Value = "1.963.41"
$ws.Range("E14").Value = "  +4.76%  "

$ws.Range("E15").Value = "  +4.46%  "

$ws.Range("D16").Value = "0.0₅8280"
$ws.Range("E16").Value = "  +2.27%  "

$ws.Range("D17").Value = "'67.79"

$ws.Range("D18").Value = "27.506.08"
$ws.Range("E18").Value = "  +5.54%  "

$ws.Range("D19").Value = "'217.82"
$ws.Range("E19").Value = "  +13.12%  "

$ws.Range("D20").Value = "'1.005"
$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("E21").Value = "  +2.56%  "

$ws.Range("E22").Value = "  +1.58%  "

$ws.Range("D23").Value = "'6.066"
$ws.Range("E23").Value = "  +2.70%  "

$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("D25").Value = "'142.57"
$ws.Range("E25").Value = "  -1.22%  "

$ws.Range("E26").Value = "  +16.14%  "

$ws.Range("E27").Value = "  +4.57%  "

$ws.Range("E28").Value = "  +2.98%  "

$ws.Range("D29").Value = "'16.50"
$ws.Range("E29").Value = "  +4.28%  "

$ws.Range("D30").Value = "'0.05486"
$ws.Range("E30").Value = "  +2.83%  "

$ws.Range("D31").Value = "'1.301"
$ws.Range("E31").Value = "  +2.58%  "

$ws.Range("D32").Value = "'3.562"
$ws.Range("E32").Value = "  +3.38%  "

$ws.Range("D33").Value = "'3.440"
$ws.Range("E33").Value = "  +3.57%  "

$ws.Range("E34").Value = "  +7.08%  "

$ws.Range("D35").Value = "'2.865"
$ws.Range("E35").Value = "  +2.99%  "

$ws.Range("D36").Value = "'0.9628"
$ws.Range("E36").Value = "  +2.41%  "

$ws.Range("D37").Value = "'2.427"
$ws.Range("E37").Value = "  +0.38%  "

$ws.Range("D38").Value = "'0.5957"
$ws.Range("E38").Value = "  +6.75%  "

$ws.Range("D39").Value = "'0.01650"
$ws.Range("E39").Value = "  +5.07%  "

$ws.Range("D40").Value = "'5.905"
$ws.Range("E40").Value = "  +2.36%  "

$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.053.38"
$ws.Range("E41").Value = "  +2.73%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.8487"
$ws.Range("E42").Value = "  +2.87%  "

$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("D44").Value = "'101.30"
$ws.Range("E44").Value = "  +0.49%  "

$ws.Range("D45").Value = "1.869.81"
$ws.Range("E45").Value = "  +4.73%  "

$ws.Range("E46").Value = "  +3.87%  "

$ws.Range("D47").Value = "'58.83"
$ws.Range("E47").Value = "  +2.75%  "

$ws.Range("D48").Value = "'0.4477"
$ws.Range("E48").Value = "  +3.85%  "

$ws.Range("D49").Value = "'8.220"
$ws.Range("E49").Value = "  +4.19%  "

$ws.Range("D50").Value = "'1.004"
$ws.Range("E50").Value = "  +0.08%  "

$ws.Range("D51").Value = "'0.05240"
$ws.Range("E51").Value = "  +2.79%  "
